$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formats (and values) from column Q into the new column R for rows 4-14,
# then overwrite with the 2020 figures from the source diff.
$ws.Range("Q4:Q14").Copy()
$ws.Range("R4:R14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("R4").Value = 2020
$ws.Range("R5").Value = 5
$ws.Range("R6").Value = 3.5
$ws.Range("R7").Value = 1.8
$ws.Range("R8").Value = 24.4
$ws.Range("R9").Value = 7.2
$ws.Range("R10").Value = 2.9
$ws.Range("R11").Value = 7.4
$ws.Range("R12").Value = 4
$ws.Range("R13").Value = 3.2
$ws.Range("R14").Value = 3.5

$ws.Range("R4:R14").Select()
